# Update existing rows 2-13 and append new rows 14-17 per natmi re-run (Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col1a1"
$ws.Cells.Item(2, 3).Value = "Cd93"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 9.108069666666667
$ws.Cells.Item(2, 8).Value = 27.324209
$ws.Cells.Item(2, 9).Value = 0.00155006418458712
$ws.Cells.Item(2, 10).Value = 0.00155006418458712
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 135.955556
$ws.Cells.Item(2, 14).Value = 407.866668
$ws.Cells.Item(2, 15).Value = 0.6947679994035034
$ws.Cells.Item(2, 16).Value = 0.6947679994035034
$ws.Cells.Item(2, 17).Value = 1238.292675618402
$ws.Cells.Item(2, 18).Value = 11144.63408056561
$ws.Cells.Item(2, 19).Value = 0.001076934992472616
$ws.Cells.Item(2, 20).Value = 0.001076934992472616

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col1a1"
$ws.Cells.Item(3, 3).Value = "Cd93"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 9.108069666666667
$ws.Cells.Item(3, 8).Value = 27.324209
$ws.Cells.Item(3, 9).Value = 0.00155006418458712
$ws.Cells.Item(3, 10).Value = 0.00155006418458712
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.449122
$ws.Cells.Item(3, 14).Value = 1.347366
$ws.Cells.Item(3, 15).Value = 0.002295129398228494
$ws.Cells.Item(3, 16).Value = 0.002295129398228494
$ws.Cells.Item(3, 17).Value = 4.090634464832667
$ws.Cells.Item(3, 18).Value = 36.815710183494
$ws.Cells.Item(3, 19).Value = 0.000003557597879186977
$ws.Cells.Item(3, 20).Value = 0.000003557597879186976

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col1a1"
$ws.Cells.Item(4, 3).Value = "Cd93"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 9.108069666666667
$ws.Cells.Item(4, 8).Value = 27.324209
$ws.Cells.Item(4, 9).Value = 0.00155006418458712
$ws.Cells.Item(4, 10).Value = 0.00155006418458712
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 56.38366533333333
$ws.Cells.Item(4, 14).Value = 169.150996
$ws.Cells.Item(4, 15).Value = 0.2881350899898248
$ws.Cells.Item(4, 16).Value = 0.2881350899898248
$ws.Cells.Item(4, 17).Value = 513.5463519180182
$ws.Cells.Item(4, 18).Value = 4621.917167262164
$ws.Cells.Item(4, 19).Value = 0.0004466278833160141
$ws.Cells.Item(4, 20).Value = 0.0004466278833160141

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Col1a1"
$ws.Cells.Item(5, 3).Value = "Cd93"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 9.108069666666667
$ws.Cells.Item(5, 8).Value = 27.324209
$ws.Cells.Item(5, 9).Value = 0.00155006418458712
$ws.Cells.Item(5, 10).Value = 0.00155006418458712
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.896484
$ws.Cells.Item(5, 14).Value = 8.689452
$ws.Cells.Item(5, 15).Value = 0.01480178120844327
$ws.Cells.Item(5, 16).Value = 0.01480178120844327
$ws.Cells.Item(5, 17).Value = 26.38137806038533
$ws.Cells.Item(5, 18).Value = 237.432402543468
$ws.Cells.Item(5, 19).Value = 0.00002294371091930257
$ws.Cells.Item(5, 20).Value = 0.00002294371091930257

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col1a1"
$ws.Cells.Item(6, 3).Value = "Cd93"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5771.873535333333
$ws.Cells.Item(6, 8).Value = 17315.620606
$ws.Cells.Item(6, 9).Value = 0.9822909543423312
$ws.Cells.Item(6, 10).Value = 0.9822909543423313
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 135.955556
$ws.Cells.Item(6, 14).Value = 407.866668
$ws.Cells.Item(6, 15).Value = 0.6947679994035034
$ws.Cells.Item(6, 16).Value = 0.6947679994035034
$ws.Cells.Item(6, 17).Value = 784718.275657929
$ws.Cells.Item(6, 18).Value = 7062464.480921361
$ws.Cells.Item(6, 19).Value = 0.6824643211805795
$ws.Cells.Item(6, 20).Value = 0.6824643211805796

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col1a1"
$ws.Cells.Item(7, 3).Value = "Cd93"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5771.873535333333
$ws.Cells.Item(7, 8).Value = 17315.620606
$ws.Cells.Item(7, 9).Value = 0.9822909543423312
$ws.Cells.Item(7, 10).Value = 0.9822909543423313
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.449122
$ws.Cells.Item(7, 14).Value = 1.347366
$ws.Cells.Item(7, 15).Value = 0.002295129398228494
$ws.Cells.Item(7, 16).Value = 0.002295129398228494
$ws.Cells.Item(7, 17).Value = 2592.275385935977
$ws.Cells.Item(7, 18).Value = 23330.4784734238
$ws.Cells.Item(7, 19).Value = 0.002254484846925008
$ws.Cells.Item(7, 20).Value = 0.002254484846925007

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Col1a1"
$ws.Cells.Item(8, 3).Value = "Cd93"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5771.873535333333
$ws.Cells.Item(8, 8).Value = 17315.620606
$ws.Cells.Item(8, 9).Value = 0.9822909543423312
$ws.Cells.Item(8, 10).Value = 0.9822909543423313
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 56.38366533333333
$ws.Cells.Item(8, 14).Value = 169.150996
$ws.Cells.Item(8, 15).Value = 0.2881350899898248
$ws.Cells.Item(8, 16).Value = 0.2881350899898248
$ws.Cells.Item(8, 17).Value = 325439.3857625581
$ws.Cells.Item(8, 18).Value = 2928954.471863023
$ws.Cells.Item(8, 19).Value = 0.2830324925256185
$ws.Cells.Item(8, 20).Value = 0.2830324925256185

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Col1a1"
$ws.Cells.Item(9, 3).Value = "Cd93"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5771.873535333333
$ws.Cells.Item(9, 8).Value = 17315.620606
$ws.Cells.Item(9, 9).Value = 0.9822909543423312
$ws.Cells.Item(9, 10).Value = 0.9822909543423313
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.896484
$ws.Cells.Item(9, 14).Value = 8.689452
$ws.Cells.Item(9, 15).Value = 0.01480178120844327
$ws.Cells.Item(9, 16).Value = 0.01480178120844327
$ws.Cells.Item(9, 17).Value = 16718.13934511643
$ws.Cells.Item(9, 18).Value = 150463.2541060479
$ws.Cells.Item(9, 19).Value = 0.01453965578920813
$ws.Cells.Item(9, 20).Value = 0.01453965578920813

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Col1a1"
$ws.Cells.Item(10, 3).Value = "Cd93"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.272029666666667
$ws.Cells.Item(10, 8).Value = 3.816089
$ws.Cells.Item(10, 9).Value = 0.0002164813950916887
$ws.Cells.Item(10, 10).Value = 0.0002164813950916887
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 135.955556
$ws.Cells.Item(10, 14).Value = 407.866668
$ws.Cells.Item(10, 15).Value = 0.6947679994035034
$ws.Cells.Item(10, 16).Value = 0.6947679994035034
$ws.Cells.Item(10, 17).Value = 172.9395005801613
$ws.Cells.Item(10, 18).Value = 1556.455505221452
$ws.Cells.Item(10, 19).Value = 0.000150404345775932
$ws.Cells.Item(10, 20).Value = 0.000150404345775932

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Col1a1"
$ws.Cells.Item(11, 3).Value = "Cd93"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.272029666666667
$ws.Cells.Item(11, 8).Value = 3.816089
$ws.Cells.Item(11, 9).Value = 0.0002164813950916887
$ws.Cells.Item(11, 10).Value = 0.0002164813950916887
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.449122
$ws.Cells.Item(11, 14).Value = 1.347366
$ws.Cells.Item(11, 15).Value = 0.002295129398228494
$ws.Cells.Item(11, 16).Value = 0.002295129398228494
$ws.Cells.Item(11, 17).Value = 0.5712965079526666
$ws.Cells.Item(11, 18).Value = 5.141668571574
$ws.Cells.Item(11, 19).Value = 0.0000004968528140444523
$ws.Cells.Item(11, 20).Value = 0.0000004968528140444523

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Col1a1"
$ws.Cells.Item(12, 3).Value = "Cd93"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.272029666666667
$ws.Cells.Item(12, 8).Value = 3.816089
$ws.Cells.Item(12, 9).Value = 0.0002164813950916887
$ws.Cells.Item(12, 10).Value = 0.0002164813950916887
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 56.38366533333333
$ws.Cells.Item(12, 14).Value = 169.150996
$ws.Cells.Item(12, 15).Value = 0.2881350899898248
$ws.Cells.Item(12, 16).Value = 0.2881350899898248
$ws.Cells.Item(12, 17).Value = 71.72169501940489
$ws.Cells.Item(12, 18).Value = 645.495255174644
$ws.Cells.Item(12, 19).Value = 0.00006237588625586654
$ws.Cells.Item(12, 20).Value = 0.00006237588625586654

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Col1a1"
$ws.Cells.Item(13, 3).Value = "Cd93"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.272029666666667
$ws.Cells.Item(13, 8).Value = 3.816089
$ws.Cells.Item(13, 9).Value = 0.0002164813950916887
$ws.Cells.Item(13, 10).Value = 0.0002164813950916887
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.896484
$ws.Cells.Item(13, 14).Value = 8.689452
$ws.Cells.Item(13, 15).Value = 0.01480178120844327
$ws.Cells.Item(13, 16).Value = 0.01480178120844327
$ws.Cells.Item(13, 17).Value = 3.684413577025333
$ws.Cells.Item(13, 18).Value = 33.159722193228
$ws.Cells.Item(13, 19).Value = 0.000003204310245845741
$ws.Cells.Item(13, 20).Value = 0.000003204310245845742

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Col1a1"
$ws.Cells.Item(14, 3).Value = "Cd93"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 93.67702500000001
$ws.Cells.Item(14, 8).Value = 281.031075
$ws.Cells.Item(14, 9).Value = 0.01594250007799006
$ws.Cells.Item(14, 10).Value = 0.01594250007799006
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 135.955556
$ws.Cells.Item(14, 14).Value = 407.866668
$ws.Cells.Item(14, 15).Value = 0.6947679994035034
$ws.Cells.Item(14, 16).Value = 0.6947679994035034
$ws.Cells.Item(14, 17).Value = 12735.9120183009
$ws.Cells.Item(14, 18).Value = 114623.2081647081
$ws.Cells.Item(14, 19).Value = 0.01107633888467535
$ws.Cells.Item(14, 20).Value = 0.01107633888467535

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Col1a1"
$ws.Cells.Item(15, 3).Value = "Cd93"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 93.67702500000001
$ws.Cells.Item(15, 8).Value = 281.031075
$ws.Cells.Item(15, 9).Value = 0.01594250007799006
$ws.Cells.Item(15, 10).Value = 0.01594250007799006
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.449122
$ws.Cells.Item(15, 14).Value = 1.347366
$ws.Cells.Item(15, 15).Value = 0.002295129398228494
$ws.Cells.Item(15, 16).Value = 0.002295129398228494
$ws.Cells.Item(15, 17).Value = 42.07241282205001
$ws.Cells.Item(15, 18).Value = 378.6517153984501
$ws.Cells.Item(15, 19).Value = 0.00003659010061025505
$ws.Cells.Item(15, 20).Value = 0.00003659010061025504

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Col1a1"
$ws.Cells.Item(16, 3).Value = "Cd93"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 93.67702500000001
$ws.Cells.Item(16, 8).Value = 281.031075
$ws.Cells.Item(16, 9).Value = 0.01594250007799006
$ws.Cells.Item(16, 10).Value = 0.01594250007799006
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 56.38366533333333
$ws.Cells.Item(16, 14).Value = 169.150996
$ws.Cells.Item(16, 15).Value = 0.2881350899898248
$ws.Cells.Item(16, 16).Value = 0.2881350899898248
$ws.Cells.Item(16, 17).Value = 5281.854027022301
$ws.Cells.Item(16, 18).Value = 47536.6862432007
$ws.Cells.Item(16, 19).Value = 0.004593593694634456
$ws.Cells.Item(16, 20).Value = 0.004593593694634455

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Col1a1"
$ws.Cells.Item(17, 3).Value = "Cd93"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 93.67702500000001
$ws.Cells.Item(17, 8).Value = 281.031075
$ws.Cells.Item(17, 9).Value = 0.01594250007799006
$ws.Cells.Item(17, 10).Value = 0.01594250007799006
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 2.896484
$ws.Cells.Item(17, 14).Value = 8.689452
$ws.Cells.Item(17, 15).Value = 0.01480178120844327
$ws.Cells.Item(17, 16).Value = 0.01480178120844327
$ws.Cells.Item(17, 17).Value = 271.3340040801
$ws.Cells.Item(17, 18).Value = 2442.0060367209
$ws.Cells.Item(17, 19).Value = 0.0002359773980699987
$ws.Cells.Item(17, 20).Value = 0.0002359773980699987

